# Add season-record columns (Wins/Losses/Ties) to the right of the
# existing data, mirroring the header style used by the other header
# cells (copy format from AC1 so the new headers get the bold/centered
# style instead of the default one).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row ---------------------------------------------------------
$headerStyleSource = $ws.Range("AC1")
$headerStyleSource.Copy()

$ws.Range("AD1").PasteSpecial(-4122)
$ws.Range("AD1").Value = "Wins"

$ws.Range("AE1").PasteSpecial(-4122)
$ws.Range("AE1").Value = "Losses"

$ws.Range("AF1").PasteSpecial(-4122)
$ws.Range("AF1").Value = "Ties"

# --- Data rows (2-49): same season record for every player -------------
$wins = 76
$losses = 85
$ties = 0

for ($row = 2; $row -le 49; $row++) {
    $ws.Cells.Item($row, 30).Value = $wins
    $ws.Cells.Item($row, 31).Value = $losses
    $ws.Cells.Item($row, 32).Value = $ties
}
